# Swap the data for rows 67<->70 and 68<->71 on the "Artfynd" sheet.
# (Columns C, I and everything from S onward are identical between the
# swapped pairs, so only A,B,D,E,F,G,H,K,M,P,Q,R need to be written.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 67 becomes the old row 70 content ---
$ws.Range("A67").Value = 111881310
$ws.Range("B67").Value = 89425
$ws.Range("D67").Value = "NT"
$ws.Range("E67").Value = 5442
$ws.Range("F67").Value = "Tallticka"
$ws.Range("G67").Value = "Porodaedalea pini"
$ws.Range("H67").Value = "(Brot.) Murrill"
$ws.Range("K67").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("P67").Value = "Valforsen, Ång"
$ws.Range("Q67").Value = 590738.9206925276
$ws.Range("R67").Value = 7040524.002523924

# --- Row 68 becomes the old row 71 content ---
$ws.Range("A68").Value = 111881322
$ws.Range("B68").Value = 56414
$ws.Range("D68").Value = "NT"
$ws.Range("E68").Value = 100049
$ws.Range("F68").Value = "Spillkråka"
$ws.Range("G68").Value = "Dryocopus martius"
$ws.Range("H68").Value = "(Linnaeus, 1758)"
$ws.Range("K68").ClearContents()
$ws.Range("M68").Value = "gammalt bo"
$ws.Range("P68").Value = "Valforsen, Ång"
$ws.Range("Q68").Value = 590615.1562677342
$ws.Range("R68").Value = 7040278.573758457

# --- Row 70 becomes the old row 67 content ---
$ws.Range("A70").Value = 111871585
$ws.Range("B70").Value = 89405
$ws.Range("D70").Value = "NT"
$ws.Range("E70").Value = 1202
$ws.Range("F70").Value = "Ullticka"
$ws.Range("G70").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H70").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("K70").Value = ""
$ws.Range("K70").NumberFormat = "General"
$ws.Range("M70").ClearContents()
$ws.Range("P70").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q70").Value = 590630.2636057099
$ws.Range("R70").Value = 7040266.929520278

# --- Row 71 becomes the old row 68 content ---
$ws.Range("A71").Value = 111870139
$ws.Range("B71").Value = 89845
$ws.Range("D71").Value = "VU"
$ws.Range("E71").Value = 1209
$ws.Range("F71").Value = "Rynkskinn"
$ws.Range("G71").Value = "Phlebia centrifuga"
$ws.Range("H71").Value = "P.Karst."
$ws.Range("K71").Value = ""
$ws.Range("K71").NumberFormat = "General"
$ws.Range("M71").ClearContents()
$ws.Range("P71").Value = "Valforsen (Valforsen), Ång"
$ws.Range("Q71").Value = 590710.4131779457
$ws.Range("R71").Value = 7040581.765558361
